{"js": "// QA feedback fix for todi_instructions.docx\n// 1) \"Sign the TODI in front of 2 witnesses and a notary public.\" ->\n//    \"Sign and date the TODI in front of 2 witnesses and a notary public.\"\n// 2) Merge the two runs around \"You cannot revoke a TODI with a will or\" /\n//    \" an unrecorded document.\" (the _GoBack bookmark that used to sit\n//    between them moves along with the cursor to the first edit spot).\n\nconst body = context.document.body;\n\n// --- Change 1: insert \" and date\" right after \"Sign\" -------------------\nconst signMatches = body.search(\"Sign the TODI in front of 2 witnesses and a notary public.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsignMatches.load(\"items\");\nawait context.sync();\n\nif (signMatches.items.length > 0) {\n  const target = signMatches.items[0];\n\n  // Narrow down to just the word \"Sign\" inside the matched sentence so we\n  // can insert immediately after it.\n  const wordMatches = target.search(\"Sign\", { matchCase: true, matchWholeWord: false });\n  wordMatches.load(\"items\");\n  await context.sync();\n\n  const signWord = wordMatches.items[0];\n  const afterInsert = signWord.insertText(\" and date\", Word.InsertLocation.after);\n  afterInsert.load(\"text\");\n  await context.sync();\n\n  // Word leaves its \"last edit\" _GoBack bookmark right where typing\n  // stopped -- recreate that at the collapsed point after \" and date\".\n  const collapsedPoint = afterInsert.getRange(Word.RangeLocation.end);\n  collapsedPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Change 2: merge \" You cannot revoke a TODI with a will or\" with\n//     \" an unrecorded document.\" into a single run, dropping the old\n//     _GoBack bookmark that used to separate them. ----------------------\nconst revokeMatches = body.search(\n  \" You cannot revoke a TODI with a will or an unrecorded document.\",\n  { matchCase: true, matchWholeWord: false }\n);\nrevokeMatches.load(\"items\");\nawait context.sync();\n\nif (revokeMatches.items.length > 0) {\n  const revokeRange = revokeMatches.items[0];\n  // Delete the whole (bookmark-spanning) range and retype it as one run\n  // so the stale _GoBack bookmark that used to sit inside it is dropped\n  // rather than merely being pushed to the end of the range.\n  revokeRange.delete();\n  await context.sync();\n\n  const anchorMatches = body.search(\"just like the TODI.\", { matchCase: true });\n  anchorMatches.load(\"items\");\n  await context.sync();\n\n  if (anchorMatches.items.length > 0) {\n    anchorMatches.items[0].insertText(\n      \" You cannot revoke a TODI with a will or an unrecorded document.\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# QA feedback fix for todi_instructions.docx\n# 1) \"Sign the TODI in front of 2 witnesses and a notary public.\" ->\n#    \"Sign and date the TODI in front of 2 witnesses and a notary public.\"\n# 2) Merge the two runs around \"You cannot revoke a TODI with a will or\" /\n#    \" an unrecorded document.\" (the _GoBack bookmark that used to sit\n#    between them moves along with the cursor to the first edit spot).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: insert \" and date\" right after \"Sign\" ---------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Sign the TODI in front of 2 witnesses and a notary public.\")\nif ($found) {\n    # Narrow to just the word \"Sign\" within the matched sentence.\n    $signRng = $d.Range($rng.Start, $rng.Start + 4)\n    $signRng.Collapse(0)   # wdCollapseEnd\n    $signRng.InsertAfter(\" and date\")\n\n    # Word drops its \"last edit\" _GoBack bookmark right where typing\n    # stopped -- recreate that at the collapsed point after \" and date\".\n    $bmPos = $signRng.End\n    $bmRng = $d.Range($bmPos, $bmPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRng) | Out-Null\n}\n\n# --- Change 2: merge \" You cannot revoke a TODI with a will or\" with\n#     \" an unrecorded document.\" into a single run, dropping the old\n#     _GoBack bookmark that used to separate them. -------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\" You cannot revoke a TODI with a will or an unrecorded document.\")\nif ($found2) {\n    # Delete the whole (bookmark-spanning) range and retype it as one run\n    # so the stale _GoBack bookmark sitting inside it is dropped instead\n    # of merely being pushed around.\n    $rng2.Text = \"\"\n    $rng2.InsertAfter(\" You cannot revoke a TODI with a will or an unrecorded document.\")\n}\n"}
